# LDLC suivi smartphones - add a new price-history snapshot column.
# A new timestamp column is inserted right before the existing "nom" /
# "url_produit" columns (DC/DD), pushing them one column to the right
# (DC->DD, DD->DE) and widening the used range from DD206 to DE206.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suivi")

# Insert a new column at DC; everything from DC onward (nom, url_produit)
# shifts right by one column (DC->DD, DD->DE).
$ws.Columns("DC:DC").Insert()

# Header row: new snapshot timestamp for the freshly inserted column.
$ws.Range("DC1").Value = "2026-02-01 14:16:35"

# Data rows 2-80 already have a running price history up through column DB,
# so the new column simply repeats that period's price (same as DB) for
# each of those rows.
for ($r = 2; $r -le 80; $r++) {
    $dbValue = $ws.Range("DB$r").Value2
    if ($null -ne $dbValue) {
        $ws.Range("DC$r").Value = $dbValue
    }
}
